$wb = $excel.ActiveWorkbook

# RQSD-RQSD sheet holds the "alternate" RPS definition; set nuclear (row 4, col B)
# to qualify (1) so nuclear is included in this alternate definition.
$ws = $wb.Worksheets.Item("RQSD-RQSD")
$ws.Range("B4").Value = 1

# Reflect the cursor having moved down to the next row (B5) after the edit,
# without disturbing which sheet tab is actually active in the workbook.
$originalActive = $wb.ActiveSheet
$ws.Activate()
$ws.Range("B5").Select()
$originalActive.Activate()
